$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 1.18
$ws.Range("K2").Value = 4.5

# Row 4
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 1.11
$ws.Range("N4").Value = 1.37
$ws.Range("O4").Value = 2.65
$ws.Range("R4").Value = 2.32
$ws.Range("S4").Value = 1.47
$ws.Range("T4").Value = 70
$ws.Range("Y4").Value = 300
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 40
$ws.Range("AC4").Value = 200
$ws.Range("AD4").Value = 8.75
$ws.Range("AE4").Value = 6.1
$ws.Range("AF4").Value = 11.5
$ws.Range("AG4").Value = 6
$ws.Range("AH4").Value = 11.25
$ws.Range("AI4").Value = 40

# Row 7
$ws.Range("G7").Value = 1.16
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 21
$ws.Range("J7").Value = 1.04
$ws.Range("K7").Value = 8.25
$ws.Range("L7").Value = 1.23
$ws.Range("M7").Value = 3.7
$ws.Range("N7").Value = 1.7
$ws.Range("O7").Value = 2.02
$ws.Range("P7").Value = 1.34
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 2.92
$ws.Range("S7").Value = 1.35
$ws.Range("T7").Value = 5.8
$ws.Range("V7").Value = 11.25
$ws.Range("W7").Value = 5.6
$ws.Range("X7").Value = 13.5
$ws.Range("Y7").Value = 60
$ws.Range("Z7").Value = 8.25
$ws.Range("AA7").Value = 14
$ws.Range("AB7").Value = 50
$ws.Range("AD7").Value = 40
$ws.Range("AE7").Value = 200
$ws.Range("AF7").Value = 75
$ws.Range("AH7").Value = 500
$ws.Range("AI7").Value = 350

# Row 8
$ws.Range("P8").Value = 1.36
$ws.Range("T8").Value = 11.5
$ws.Range("U8").Value = 22
$ws.Range("AC8").Value = 60
$ws.Range("AD8").Value = 8.5
$ws.Range("AE8").Value = 10.75
$ws.Range("AF8").Value = 8.75

# Row 9
$ws.Range("H9").Value = 3.05
$ws.Range("I9").Value = 3.1
$ws.Range("K9").Value = 6.2
$ws.Range("M9").Value = 2.65
$ws.Range("N9").Value = 2.3
$ws.Range("R9").Value = 1.93
$ws.Range("S9").Value = 1.78
$ws.Range("T9").Value = 6.4
$ws.Range("U9").Value = 11
$ws.Range("W9").Value = 26
$ws.Range("Z9").Value = 6.2
$ws.Range("AA9").Value = 6.2
$ws.Range("AB9").Value = 17.5
$ws.Range("AD9").Value = 7.9
$ws.Range("AE9").Value = 16
$ws.Range("AG9").Value = 45
$ws.Range("AJ9").Value = 1000

# Row 12
$ws.Range("N12").Value = 2.05
$ws.Range("O12").Value = 1.75
$ws.Range("T12").Value = 6
$ws.Range("AB12").Value = 19

# Row 19
$ws.Range("AD19").Value = 13.5
$ws.Range("AE19").Value = 26
$ws.Range("AI19").Value = 40

# Row 20
$ws.Range("U20").Value = 8.25
$ws.Range("X20").Value = 12.5
$ws.Range("Z20").Value = 11.25

# Row 21
$ws.Range("G21").Value = 1.8
$ws.Range("H21").Value = 3.45
$ws.Range("I21").Value = 4.05
$ws.Range("T21").Value = 7.3
$ws.Range("Z21").Value = 10
$ws.Range("AB21").Value = 14.5
$ws.Range("AC21").Value = 65
$ws.Range("AD21").Value = 11.25
$ws.Range("AG21").Value = 65

# Row 22
$ws.Range("G22").Value = 1.83
$ws.Range("H22").Value = 3.5
$ws.Range("L22").Value = 1.31
$ws.Range("M22").Value = 2.87
$ws.Range("N22").Value = 1.91
$ws.Range("O22").Value = 1.7
$ws.Range("R22").Value = 1.83
$ws.Range("S22").Value = 1.78
$ws.Range("T22").Value = 6.6
$ws.Range("V22").Value = 8.5
$ws.Range("X22").Value = 15.5
$ws.Range("Y22").Value = 30
$ws.Range("Z22").Value = 9.25
$ws.Range("AB22").Value = 16.5
$ws.Range("AC22").Value = 80
$ws.Range("AD22").Value = 10.5
$ws.Range("AJ22").Value = 700

# Row 23
$ws.Range("AD23").Value = 6.7
$ws.Range("AE23").Value = 6.9

# Row 24
$ws.Range("G24").Value = 2.72
$ws.Range("H24").Value = 3.15
$ws.Range("I24").Value = 2.35
$ws.Range("N24").Value = 1.91
$ws.Range("O24").Value = 1.7
$ws.Range("T24").Value = 7.4
$ws.Range("U24").Value = 11.5
$ws.Range("V24").Value = 8.5
$ws.Range("W24").Value = 25
$ws.Range("X24").Value = 18.5
$ws.Range("Y24").Value = 25
$ws.Range("AA24").Value = 5.4
$ws.Range("AB24").Value = 11.5
$ws.Range("AC24").Value = 50
$ws.Range("AD24").Value = 6.6
$ws.Range("AE24").Value = 9.5
$ws.Range("AF24").Value = 7.9
$ws.Range("AG24").Value = 19
$ws.Range("AH24").Value = 16
$ws.Range("AJ24").Value = 350

# Row 25
$ws.Range("T25").Value = 6.2
$ws.Range("Y25").Value = 24
$ws.Range("Z25").Value = 13.5
$ws.Range("AB25").Value = 20
$ws.Range("AD25").Value = 19.5
$ws.Range("AI25").Value = 75

# Row 28
$ws.Range("T28").Value = 6.3
$ws.Range("X28").Value = 21
$ws.Range("Y28").Value = 37
$ws.Range("Z28").Value = 7.2
$ws.Range("AD28").Value = 7.7
$ws.Range("AE28").Value = 15
$ws.Range("AF28").Value = 11.75
$ws.Range("AI28").Value = 50

# Row 29
$ws.Range("M29").Value = 2.4
$ws.Range("N29").Value = 2.27
$ws.Range("W29").Value = 23
$ws.Range("Y29").Value = 40
$ws.Range("AB29").Value = 17
$ws.Range("AD29").Value = 7.5
$ws.Range("AF29").Value = 11.75
$ws.Range("AG29").Value = 40
$ws.Range("AI29").Value = 50

# Row 32
$ws.Range("G32").Value = 2.5
$ws.Range("H32").Value = 3.5
$ws.Range("I32").Value = 2.6
$ws.Range("L32").Value = 1.18
$ws.Range("M32").Value = 4.5
$ws.Range("N32").Value = 1.62
$ws.Range("O32").Value = 2.25
$ws.Range("P32").Value = 1.29
$ws.Range("Q32").Value = 3.5
$ws.Range("R32").Value = 1.5
$ws.Range("S32").Value = 2.5
$ws.Range("T32").Value = 12
$ws.Range("AA32").Value = 7
$ws.Range("AB32").Value = 11
$ws.Range("AC32").Value = 34
$ws.Range("AD32").Value = 12
$ws.Range("AE32").Value = 15
$ws.Range("AF32").Value = 10
$ws.Range("AG32").Value = 26
$ws.Range("AJ32").Value = 101

# Row 33
$ws.Range("L33").Value = 1.17
$ws.Range("M33").Value = 5
$ws.Range("N33").Value = 1.6
$ws.Range("O33").Value = 2.3

# Row 35
$ws.Range("K35").Value = 23

